# Insert a new "OrthogonalImage" worksheet right after "ImageMask" (and
# before "ChannelSeries"), matching the schema table used by the other
# "Image"-like sheets plus source_image / source_roi / axis columns, and
# an "axis" column restricted to XY/XZ/YZ via a list data validation.

$wb = $excel.ActiveWorkbook

$imageMask = $wb.Worksheets.Item("ImageMask")
$newSheet = $wb.Worksheets.Add($null, $imageMask)
$newSheet.Name = "OrthogonalImage"

$headers = @(
    "source_image",
    "source_roi",
    "axis",
    "voxel_size_x_micron",
    "voxel_size_y_micron",
    "voxel_size_z_micron",
    "shape_x",
    "shape_y",
    "shape_z",
    "shape_c",
    "shape_t",
    "time_series",
    "channel_series",
    "acquisition_datetime",
    "source_images",
    "array_data",
    "data_reference",
    "linked_references",
    "name",
    "description"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# "axis" (column C) is a dropdown restricted to XY / XZ / YZ, for every
# row below the header.
$axisRange = $newSheet.Range("C2:C1048576")
$axisRange.Validation.Add(3, 1, 1, '"XY,XZ,YZ"')
$axisRange.Validation.ShowInput = $false
$axisRange.Validation.ShowError = $false
